$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "producto" in K1, styled like the other header cells (C1..J1)
# but bold: copy the existing header format then flip Bold on.
$ws.Range("K1").Value = "producto"
$ws.Range("C1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Font.Bold = $true

# New data column K2:K6 with value 3
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 3

# Update the active selection shown in the sheet view
$ws.Range("L9").Select()
